# Adds support for a new "text" function category (`outputToCloud(resource)`
# moves into `base`, and a brand-new `text` category holding
# `spellCheck(var,profile,text)`) to the hidden '#system' lookup sheet that
# backs several of the workbook's data-validation dropdown-driven named
# ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1) "target" category list (column A) gains a new category: "text",
#    inserted alphabetically between "step" and "web" (row 25), shifting the
#    remaining categories down one row.
$ws.Range("A25").Insert(-4121)
$ws.Range("A25").Value = "text"

# 2) "base" function list (column E) gains a new function:
#    "outputToCloud(resource)", inserted alphabetically between
#    "macro(file,sheet,name)" and "prependText(var,prependWith)" (row 22),
#    shifting the remaining functions down one row.
$ws.Range("E22").Insert(-4121)
$ws.Range("E22").Value = "outputToCloud(resource)"

# 3) A brand-new column is inserted for the "text" category's function list,
#    immediately to the left of the existing "web" column (column Y),
#    shifting "web" through "xml" one column to the right.
$ws.Range("Y1:Y129").Insert(-4161)
$ws.Range("Y1").Value = "text"
$ws.Range("Y2").Value = "spellCheck(var,profile,text)"

# 4) Update the defined names that pointed at the shifted ranges.
$wb.Names.Item("base").RefersTo = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AE`$2:`$AE`$27"

# 5) Register the new "text" category as its own named range.
$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")
